# Implemented Inventory Extract Report logic
# - Remove the "OrchestratorQueueName" row from Settings
# - Rework the Assets sheet to add an "OrchestratorAssetFolder" column and
#   replace/extend the asset rows with the JDE-related assets
# - Add a new "Credentials" sheet holding the JDE credential asset

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Settings sheet - drop the Orchestrator Queue Name setting (row 10)
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Activate()
$wsSettings.Range("A10:C10").ClearContents()
$wsSettings.Rows.Item(10).Select()

# ---------------------------------------------------------------------
# 2) Assets sheet - add OrchestratorAssetFolder column and new asset rows
# ---------------------------------------------------------------------
$wsAssets = $wb.Worksheets.Item("Assets")
$wsAssets.Activate()

# Header row
$wsAssets.Range("A1").Value = "Name"
$wsAssets.Range("B1").Value = "Value"
$wsAssets.Range("C1").Value = "OrchestratorAssetFolder"
$wsAssets.Range("A1:C1").Font.Bold = $true
$wsAssets.Range("A1:C1").Font.Size = 14
$wsAssets.Range("A1:C1").Font.Name = "Calibri"
$wsAssets.Range("A1:C1").Font.Color = 0x000000
$wsAssets.Rows.Item(1).RowHeight = 18.5

# Existing assets keep their name/value, gain a folder of "Corporate"
$wsAssets.Range("A2").Value = "796_GetEmail_Count"
$wsAssets.Range("B2").Value = "796_GetEmail_Count"
$wsAssets.Range("C2").Value = "Corporate"

$wsAssets.Range("A3").Value = "796_Download_Reports_Path"
$wsAssets.Range("B3").Value = "796_Download_Reports_Path"
$wsAssets.Range("C3").Value = "Corporate"

# New assets for the JDE integration
$wsAssets.Range("A4").Value = "796_JDE_Url"
$wsAssets.Range("B4").Value = "796_JDE_Url"
$wsAssets.Range("C4").Value = "Corporate"

$wsAssets.Range("A5").Value = "796_Role"
$wsAssets.Range("B5").Value = "796_Role"
$wsAssets.Range("C5").Value = "Corporate"

$wsAssets.Range("A6").Value = "796_Queue_Name"
$wsAssets.Range("B6").Value = "796_Queue_Name"
$wsAssets.Range("C6").Value = "Corporate"

$wsAssets.Range("A5").Select()

# ---------------------------------------------------------------------
# 3) Credentials sheet (new) - holds the JDE credential asset
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsCred = $wb.Worksheets.Add($null, $lastSheet)
$wsCred.Name = "Credentials"

$wsCred.Columns.Item(1).ColumnWidth = 46.6328125
$wsCred.Columns.Item(2).ColumnWidth = 33.08984375
$wsCred.Columns.Item(3).ColumnWidth = 27.453125

$wsCred.Range("A1").Value = "Name"
$wsCred.Range("B1").Value = "Value"
$wsCred.Range("C1").Value = "OrchestratorAssetFolder"
$wsCred.Range("A1:C1").Font.Bold = $true
$wsCred.Range("A1:C1").Font.Size = 14
$wsCred.Range("A1:C1").Font.Name = "Calibri"
$wsCred.Range("A1:C1").Font.Color = 0x000000
$wsCred.Rows.Item(1).RowHeight = 18.5

$wsCred.Range("A2").Value = "796_JDE_Credentials"
$wsCred.Range("B2").Value = "796_JDE_Credentials"
$wsCred.Range("A2:B2").Font.Size = 7
$wsCred.Range("A2:B2").Font.Name = "Noto Sans"
$wsCred.Range("A2:B2").Font.Color = 0x554E46
$wsCred.Range("C2").Value = "Corporate"

$wsCred.Range("A2").Select()
$wsCred.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 4) Leave the Assets sheet as the active tab (matches the saved workbook)
# ---------------------------------------------------------------------
$wsAssets.Activate()
